$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 616.3333
$ws.Range("I33").Value = 221.1
$ws.Range("K33").Value = 221.1
$ws.Range("M33").Value = 7.900000000000006

$ws.Range("H129").Value = 672.35297
$ws.Range("I129").Value = 502.27274
$ws.Range("J129").Value = 984.1667
$ws.Range("K129").Value = 1506.81822
$ws.Range("L129").Value = 2952.5001
$ws.Range("M129").Value = 3493.18178
$ws.Range("N129").Value = -12952.5001

$ws.Range("H131").Value = 6590.609
$ws.Range("I131").Value = 638.93335
$ws.Range("J131").Value = 17750
$ws.Range("K131").Value = 1916.80005
$ws.Range("L131").Value = 53250
$ws.Range("M131").Value = 3123.19995
$ws.Range("N131").Value = -63330

$ws.Range("H132").Value = 156171.72
$ws.Range("I132").Value = 2167.261
$ws.Range("J132").Value = 529024.6
$ws.Range("K132").Value = 6501.782999999999
$ws.Range("L132").Value = 1587073.8
$ws.Range("M132").Value = -3971.782999999999
$ws.Range("N132").Value = -1592133.8

$ws.Range("H135").Value = 9260195
$ws.Range("I135").Value = 285.40625
$ws.Range("J135").Value = 22729154
$ws.Range("K135").Value = 2568.65625
$ws.Range("L135").Value = 204562386
$ws.Range("M135").Value = -33.65625
$ws.Range("N135").Value = -204567456

$ws.Range("H137").Value = 45261
$ws.Range("I137").Value = 92256.55
$ws.Range("J137").Value = 8335.929
$ws.Range("K137").Value = 276769.65
$ws.Range("L137").Value = 25007.787
$ws.Range("M137").Value = -274219.65
$ws.Range("N137").Value = -30107.787

$ws.Range("H138").Value = 1673.62
$ws.Range("I138").Value = 850.0465
$ws.Range("J138").Value = 2294.9124
$ws.Range("K138").Value = 2550.1395
$ws.Range("L138").Value = 6884.7372
$ws.Range("M138").Value = 2589.8605
$ws.Range("N138").Value = -17164.7372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 19403.709
$ws.Range("I74").Value = 29791.686
$ws.Range("J74").Value = 1224.75
$ws.Range("K74").Value = 29791.686
$ws.Range("L74").Value = 1224.75
$ws.Range("M74").Value = -28917.686
$ws.Range("N74").Value = -2972.75

$ws.Range("H77").Value = 19403.709
$ws.Range("I77").Value = 29791.686
$ws.Range("J77").Value = 1224.75
$ws.Range("K77").Value = 148958.43
$ws.Range("L77").Value = 6123.75
$ws.Range("M77").Value = -144590.43
$ws.Range("N77").Value = -14859.75

$ws.Range("H122").Value = 1156
$ws.Range("I122").Value = 1156
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3468
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1018
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1056.3654
$ws.Range("I58").Value = 760.3214
$ws.Range("K58").Value = 760.3214
$ws.Range("M58").Value = -557.3214

$ws.Range("H94").Value = 1559.1923
$ws.Range("I94").Value = 1831.6
$ws.Range("J94").Value = 1388.9375
$ws.Range("K94").Value = 1831.6
$ws.Range("L94").Value = 1388.9375
$ws.Range("M94").Value = -1380.6
$ws.Range("N94").Value = -2290.9375

$ws.Range("H132").Value = 940.04
$ws.Range("I132").Value = 716.02856
$ws.Range("K132").Value = 2148.08568
$ws.Range("M132").Value = 381.9143199999999

$ws.Range("H134").Value = 879.66174
$ws.Range("I134").Value = 809.1064
$ws.Range("J134").Value = 1037.5714
$ws.Range("K134").Value = 2427.3192
$ws.Range("L134").Value = 3112.7142
$ws.Range("M134").Value = 107.6808000000001
$ws.Range("N134").Value = -8182.7142

$ws.Range("H136").Value = 1056.3654
$ws.Range("I136").Value = 760.3214
$ws.Range("K136").Value = 2280.9642
$ws.Range("M136").Value = 269.0357999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3475.543
$ws.Range("I5").Value = 722.05
$ws.Range("J5").Value = 7146.8667
$ws.Range("K5").Value = 2166.15
$ws.Range("L5").Value = 21440.6001
$ws.Range("M5").Value = -2054.15
$ws.Range("N5").Value = -21664.6001

$ws.Range("H113").Value = 477.66666
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 484.72726
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 1454.18178
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -5794.18178

$ws.Range("H122").Value = 402.45715
$ws.Range("I122").Value = 359.5
$ws.Range("J122").Value = 438.6316
$ws.Range("K122").Value = 3235.5
$ws.Range("L122").Value = 3947.6844
$ws.Range("M122").Value = -785.5
$ws.Range("N122").Value = -8847.6844

$ws.Range("H135").Value = 3475.543
$ws.Range("I135").Value = 722.05
$ws.Range("J135").Value = 7146.8667
$ws.Range("K135").Value = 6498.45
$ws.Range("L135").Value = 64321.8003
$ws.Range("M135").Value = -3963.45
$ws.Range("N135").Value = -69391.8003

$ws.Range("H137").Value = 8688997
$ws.Range("I137").Value = 23810282
$ws.Range("J137").Value = 4455037
$ws.Range("K137").Value = 71430846
$ws.Range("L137").Value = 13365111
$ws.Range("M137").Value = -71425746
$ws.Range("N137").Value = -13375311

$ws.Range("H138").Value = 14495484
$ws.Range("I138").Value = 1610
$ws.Range("J138").Value = 27781534
$ws.Range("K138").Value = 4830
$ws.Range("L138").Value = 83344602
$ws.Range("M138").Value = 310
$ws.Range("N138").Value = -83354882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1394.26
$ws.Range("I126").Value = 1272.0278
$ws.Range("J126").Value = 1708.5714
$ws.Range("K126").Value = 3816.0834
$ws.Range("L126").Value = 5125.7142
$ws.Range("M126").Value = -1346.0834
$ws.Range("N126").Value = -10065.7142

$ws.Range("H132").Value = 20980.883
$ws.Range("I132").Value = 1219.7179
$ws.Range("K132").Value = 3659.1537
$ws.Range("M132").Value = -1129.1537

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2716.1353
$ws.Range("I7").Value = 1689.56
$ws.Range("J7").Value = 4854.8335
$ws.Range("K7").Value = 1689.56
$ws.Range("L7").Value = 4854.8335
$ws.Range("M7").Value = -1577.56
$ws.Range("N7").Value = -5078.8335

$ws.Range("H40").Value = 37892.25
$ws.Range("I40").Value = 1350.5294
$ws.Range("J40").Value = 94365.82000000001
$ws.Range("K40").Value = 1350.5294
$ws.Range("L40").Value = 94365.82000000001
$ws.Range("M40").Value = -1214.5294
$ws.Range("N40").Value = -94637.82000000001

$ws.Range("H61").Value = 2824.9167
$ws.Range("I61").Value = 2599.8333
$ws.Range("K61").Value = 2599.8333
$ws.Range("M61").Value = -2397.8333

$ws.Range("H113").Value = 2824.9167
$ws.Range("I113").Value = 2599.8333
$ws.Range("K113").Value = 2599.8333
$ws.Range("M113").Value = -429.8332999999998

$ws.Range("H122").Value = 3247.125
$ws.Range("I122").Value = 3322.6155
$ws.Range("J122").Value = 2920
$ws.Range("K122").Value = 9967.8465
$ws.Range("L122").Value = 8760
$ws.Range("M122").Value = -7517.8465
$ws.Range("N122").Value = -13660

$ws.Range("H126").Value = 2716.1353
$ws.Range("I126").Value = 1689.56
$ws.Range("J126").Value = 4854.8335
$ws.Range("K126").Value = 5068.68
$ws.Range("L126").Value = 14564.5005
$ws.Range("M126").Value = -2598.68
$ws.Range("N126").Value = -19504.5005

$ws.Range("H132").Value = 193072.95
$ws.Range("I132").Value = 40815.137
$ws.Range("J132").Value = 840168.7
$ws.Range("K132").Value = 122445.411
$ws.Range("L132").Value = 2520506.1
$ws.Range("M132").Value = -119915.411
$ws.Range("N132").Value = -2525566.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9572.727999999999
$ws.Range("I122").Value = 9450
$ws.Range("K122").Value = 28350
$ws.Range("M122").Value = -25900

$ws.Range("H126").Value = 1336.1428
$ws.Range("I126").Value = 1030.6
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 3091.8
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -621.7999999999997
$ws.Range("N126").Value = -11240

$ws.Range("H132").Value = 3417.9756
$ws.Range("I132").Value = 1023.5
$ws.Range("J132").Value = 6190.5264
$ws.Range("K132").Value = 3070.5
$ws.Range("L132").Value = 18571.5792
$ws.Range("M132").Value = -540.5
$ws.Range("N132").Value = -23631.5792
